$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 1.033366029227541
$ws.Range("E2").Value = 1.041296853725589
$ws.Range("F2").Value = 1.049333393264742
$ws.Range("J2").Value = 1.038491151223592
$ws.Range("L2").Value = 1.044076768678686
$ws.Range("M2").Value = 1.052090743175996
$ws.Range("N2").Value = 1.016359438818662
$ws.Range("C3").Value = 1.035456500123056
$ws.Range("E3").Value = 1.043210675826196
$ws.Range("F3").Value = 1.051452051245787
$ws.Range("J3").Value = 1.040218527726439
$ws.Range("L3").Value = 1.045798295184488
$ws.Range("M3").Value = 1.054018255724803
$ws.Range("N3").Value = 1.016972109842518
$ws.Range("C4").Value = 1.036804338530668
$ws.Range("E4").Value = 1.044444645201333
$ws.Range("F4").Value = 1.052818661432288
$ws.Range("J4").Value = 1.041331346126494
$ws.Range("L4").Value = 1.046907440178717
$ws.Range("M4").Value = 1.055260824523871
$ws.Range("N4").Value = 1.017365833892181
$ws.Range("C5").Value = 1.037369841693725
$ws.Range("E5").Value = 1.044962378131683
$ws.Range("F5").Value = 1.053392182731226
$ws.Range("J5").Value = 1.041798024285772
$ws.Range("L5").Value = 1.047372600427122
$ws.Range("M5").Value = 1.055782111355798
$ws.Range("N5").Value = 1.017530713907021
$ws.Range("C6").Value = 1.037464726867819
$ws.Range("E6").Value = 1.045049248248011
$ws.Range("F6").Value = 1.05348842157114
$ws.Range("J6").Value = 1.041876314838629
$ws.Range("L6").Value = 1.047450637651735
$ws.Range("M6").Value = 1.055869574615348
$ws.Range("N6").Value = 1.01755836060787
$ws.Range("C7").Value = 1.036811899204402
$ws.Range("E7").Value = 1.0444515671777
$ws.Range("F7").Value = 1.052826328749198
$ws.Range("J7").Value = 1.0413375863911
$ws.Range("L7").Value = 1.046913660057775
$ws.Range("M7").Value = 1.055267794227486
$ws.Range("N7").Value = 1.017368039536576
$ws.Range("C8").Value = 1.034073532951003
$ws.Range("E8").Value = 1.041944565930772
$ws.Range("F8").Value = 1.050050311958026
$ws.Range("J8").Value = 1.039075957746021
$ws.Range("L8").Value = 1.04465957474345
$ws.Range("M8").Value = 1.05274313626363
$ws.Range("N8").Value = 1.016567061135038
$ws.Range("C9").Value = 1.029209801128962
$ws.Range("E9").Value = 1.037492025308603
$ws.Range("F9").Value = 1.04512436522637
$ws.Range("J9").Value = 1.035051958522593
$ws.Range("L9").Value = 1.040649743577525
$ws.Range("M9").Value = 1.048257462646568
$ws.Range("N9").Value = 1.015134458890293
$ws.Range("C10").Value = 1.025939617291361
$ws.Range("E10").Value = 1.034498556500014
$ws.Range("F10").Value = 1.041815550827718
$ws.Range("J10").Value = 1.032341682945348
$ws.Range("L10").Value = 1.037949528146232
$ws.Range("M10").Value = 1.045240498997321
$ws.Range("N10").Value = 1.014164613834183
$ws.Range("C11").Value = 1.024516597594766
$ws.Range("E11").Value = 1.033196024245223
$ws.Range("F11").Value = 1.04037649423181
$ws.Range("J11").Value = 1.031161200832586
$ws.Range("L11").Value = 1.036773553826061
$ws.Range("M11").Value = 1.043927448861384
$ws.Range("N11").Value = 1.013741033680938
$ws.Range("C12").Value = 1.023986934918853
$ws.Range("E12").Value = 1.032711220915081
$ws.Range("F12").Value = 1.039840979805728
$ws.Range("J12").Value = 1.030721647574991
$ws.Range("L12").Value = 1.036335698255755
$ws.Range("M12").Value = 1.043438686032645
$ws.Range("N12").Value = 1.013583141053325
$ws.Range("C13").Value = 1.024100599274608
$ws.Range("E13").Value = 1.032815257965508
$ws.Range("F13").Value = 1.039955894605356
$ws.Range("J13").Value = 1.030815982143405
$ws.Range("L13").Value = 1.036429667602275
$ws.Range("M13").Value = 1.043543574808477
$ws.Range("N13").Value = 1.013617034894785
$ws.Range("C14").Value = 1.024472837959735
$ws.Range("E14").Value = 1.033155970530062
$ws.Range("F14").Value = 1.040332248738053
$ws.Range("J14").Value = 1.031124889223288
$ws.Range("L14").Value = 1.036737382070019
$ws.Range("M14").Value = 1.043887068933161
$ws.Range("N14").Value = 1.013727993633358
$ws.Range("C15").Value = 1.024702040976964
$ws.Range("E15").Value = 1.033365763308621
$ws.Range("F15").Value = 1.040564001381582
$ws.Range("J15").Value = 1.031315074403884
$ws.Range("L15").Value = 1.036926835544605
$ws.Range("M15").Value = 1.044098568512556
$ws.Range("N15").Value = 1.013796284965319
$ws.Range("C16").Value = 1.026033907629089
$ws.Range("E16").Value = 1.034584864946815
$ws.Range("F16").Value = 1.04191092014161
$ws.Range("J16").Value = 1.032419879352861
$ws.Range("L16").Value = 1.038027428624988
$ws.Range("M16").Value = 1.045327498144159
$ws.Range("N16").Value = 1.014192648055195
$ws.Range("C17").Value = 1.02686745018941
$ws.Range("E17").Value = 1.035347855050773
$ws.Range("F17").Value = 1.042754090270336
$ws.Range("J17").Value = 1.03311102144153
$ws.Range("L17").Value = 1.038715969682559
$ws.Range("M17").Value = 1.046096560209594
$ws.Range("N17").Value = 1.014440296640642
$ws.Range("C18").Value = 1.027352967761147
$ws.Range("E18").Value = 1.035792284708298
$ws.Range("F18").Value = 1.043245289960371
$ws.Range("J18").Value = 1.033513487843186
$ws.Range("L18").Value = 1.039116933547571
$ws.Range("M18").Value = 1.046544498368184
$ws.Range("N18").Value = 1.014584396272636
$ws.Range("C19").Value = 1.027518403541066
$ws.Range("E19").Value = 1.035943721312831
$ws.Range("F19").Value = 1.043412674438442
$ws.Range("J19").Value = 1.033650606598615
$ws.Range("L19").Value = 1.039253542438254
$ws.Range("M19").Value = 1.04669712568265
$ws.Range("N19").Value = 1.014633471543508
$ws.Range("C20").Value = 1.026778088910315
$ws.Range("E20").Value = 1.03526605671531
$ws.Range("F20").Value = 1.042663689154085
$ws.Range("J20").Value = 1.033036937450379
$ws.Range("L20").Value = 1.038642163243584
$ws.Range("M20").Value = 1.046014113847934
$ws.Range("N20").Value = 1.014413762527972
$ws.Range("C21").Value = 1.024363253327549
$ws.Range("E21").Value = 1.033055666620427
$ws.Range("F21").Value = 1.040221449297814
$ws.Range("J21").Value = 1.031033953564367
$ws.Range("L21").Value = 1.036646796958576
$ws.Range("M21").Value = 1.043785947385486
$ws.Range("N21").Value = 1.013695334487559
$ws.Range("C22").Value = 1.022838626403462
$ws.Range("E22").Value = 1.031660192195375
$ws.Range("F22").Value = 1.038680200795731
$ws.Range("J22").Value = 1.029768393883549
$ws.Range("L22").Value = 1.035386161784353
$ws.Range("M22").Value = 1.042378990933126
$ws.Range("N22").Value = 1.013240407459075
$ws.Range("C23").Value = 1.023647471882665
$ws.Range("E23").Value = 1.032400512166617
$ws.Range("F23").Value = 1.039497799546495
$ws.Range("J23").Value = 1.03043988950466
$ws.Range("L23").Value = 1.036055033862708
$ws.Range("M23").Value = 1.043125427002245
$ws.Range("N23").Value = 1.013481881958445
$ws.Range("C24").Value = 1.02681846950613
$ws.Range("E24").Value = 1.03530301974448
$ws.Range("F24").Value = 1.042704539402102
$ws.Range("J24").Value = 1.033070414871012
$ws.Range("L24").Value = 1.03867551520525
$ws.Range("M24").Value = 1.046051369794694
$ws.Range("N24").Value = 1.01442575323124
$ws.Range("C25").Value = 1.030471933579147
$ws.Range("E25").Value = 1.03864741647204
$ws.Range("F25").Value = 1.046402085227728
$ws.Range("J25").Value = 1.036097008098111
$ws.Range("L25").Value = 1.041691023941353
$ws.Range("M25").Value = 1.049421663745579
$ws.Range("N25").Value = 1.015507385343154
